$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 79, shifting rows 79:169 down to 80:170
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with the latest week's data
$ws.Range("A79").Value = 7
$ws.Range("B79").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C79").Value = "Ñuble"
$ws.Range("D79").Value = 44494
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100108
$ws.Range("H79").Value = "Tropicales y subtropicales"
$ws.Range("I79").Value = 100108005
$ws.Range("J79").Value = "Piña"
$ws.Range("K79").Value = "Caramelo"
$ws.Range("L79").Value = "Segunda"
$ws.Range("M79").Value = 120
$ws.Range("N79").Value = 20000
$ws.Range("O79").Value = 21000
$ws.Range("P79").Value = 20500
$ws.Range("Q79").Value = "$/caja 14 unidades"
$ws.Range("R79").Value = "Ecuador"
$ws.Range("S79").Value = 1464
$ws.Range("T79").Value = 14

# Match the date formatting used by the rest of column D
$ws.Range("D79").NumberFormat = $ws.Range("D80").NumberFormat
